$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.815.00'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.639.76'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.10'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.72'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.865.36'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').Value = '1.638.37'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.08'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '25.845.51'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.80'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.98'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.34'
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.84'
$ws.Range('E24').Value = '  +4.61%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.06'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  +1.91%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0494'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.908'
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('D37').Value = '1.133.05'
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.53'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.84'
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.807'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('D45').Value = '1.775.23'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').Value = '  +3.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.38'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.416'
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('E49').Value = '  +5.61%  '
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0956'
$ws.Range('E51').Value = '  +1.78%  '
